# Applies the "Add files via upload" commit:
#   1. Inserts a new worksheet "Leaderboard Utveckling" right after "Leaderboard",
#      containing the per-round ("Deltävling 1/2/3") point breakdown for every player.
#   2. Adds a new player ("Dempa") as row 12 of the "Leaderboard" sheet.
#   3. Renames the player in row 4 of "Leaderboard" from "Johan" to "Rantzow"
#      (the row's photo in column A is left untouched).
#
# NOTE: the order in which new cell values are assigned below is deliberate -
# it reproduces the order new entries were appended to the shared-strings
# table in the authored workbook.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1) New sheet "Leaderboard Utveckling", placed after "Leaderboard"
# ---------------------------------------------------------------------------
$wsLeaderboard = $wb.Worksheets.Item("Leaderboard")
$wsDev = $wb.Worksheets.Add($null, $wsLeaderboard)
$wsDev.Name = "Leaderboard Utveckling"

$players = @("Axel", "Jojo", "Johan", "Sebbe", "Alvin", "Löken", "Benne", "Crille", "Vigge", "Frasse")

$points1 = @(10, 6, 8, 4, 3, 2, 1, 1, 1, 0)
$points2 = @(21, 16, 10, 5, 7, 2, 1, 11, 15, 20)
$points3 = @(30, 20, 13, 10, 8, 8, 11, 13, 22, 28)

$wsDev.Cells.Item(1, 1).Value = "Spelare"
$wsDev.Cells.Item(1, 2).Value = "Poäng"
$wsDev.Cells.Item(1, 3).Value = "Deltävling"

$row = 2
for ($i = 0; $i -lt $players.Length; $i++) {
    $wsDev.Cells.Item($row, 1).Value = $players[$i]
    $wsDev.Cells.Item($row, 2).Value = $points1[$i]
    $wsDev.Cells.Item($row, 3).Value = "Deltävling 1"
    $row = $row + 1
}
for ($i = 0; $i -lt $players.Length; $i++) {
    $wsDev.Cells.Item($row, 1).Value = $players[$i]
    $wsDev.Cells.Item($row, 2).Value = $points2[$i]
    $wsDev.Cells.Item($row, 3).Value = "Deltävling 2"
    $row = $row + 1
}
for ($i = 0; $i -lt $players.Length; $i++) {
    $wsDev.Cells.Item($row, 1).Value = $players[$i]
    $wsDev.Cells.Item($row, 2).Value = $points3[$i]
    $wsDev.Cells.Item($row, 3).Value = "Deltävling 3"
    $row = $row + 1
}

$wsDev.Columns.Item(3).ColumnWidth = 11.5
[void]$wsDev.Range("I17").Select()

# ---------------------------------------------------------------------------
# 2) New player row (Dempa) on the "Leaderboard" sheet
# ---------------------------------------------------------------------------
$dempaPhoto = "https://scontent-arn2-1.xx.fbcdn.net/v/t31.18172-8/289992_327758097305116_1434745224_o.jpg?_nc_cat=109&ccb=1-7&_nc_sid=5f2048&_nc_ohc=rd6JEllRr2UAX89sqwh&_nc_ht=scontent-arn2-1.xx&cb_e2o_trans=q&oh=00_AfAIEQdFV33zK2qkeNYOLw7erpCeQ2zMPbZjMsrQaImenQ&oe=663248F7"

$wsLeaderboard.Range("B12").Value = "Dempa"
$wsLeaderboard.Range("A12").Value = $dempaPhoto
$wsLeaderboard.Range("C12").Value = 0
$wsLeaderboard.Range("D12").Value = 0
$wsLeaderboard.Range("E12").Value = 0

# ---------------------------------------------------------------------------
# 3) Rename player in row 4 of "Leaderboard": Johan -> Rantzow
# ---------------------------------------------------------------------------
$wsLeaderboard.Range("B4").Value = "Rantzow"

[void]$wsLeaderboard.Range("C17").Select()
